$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 00:40"

# Update country statistics cells
# Row 4
$ws.Range("B4").Value = 8443499
$ws.Range("C4").Value = 45398
$ws.Range("D4").Value = 5491461
$ws.Range("E4").Value = 2726887
$ws.Range("G4").Value = 377
$ws.Range("H4").Value = 225151

# Row 6
$ws.Range("B6").Value = 5250727
$ws.Range("C6").Value = 15383
$ws.Range("D6").Value = 4681659
$ws.Range("E6").Value = 414892
$ws.Range("G6").Value = 271
$ws.Range("H6").Value = 154176

# Row 10
$ws.Range("B10").Value = 965883
$ws.Range("C10").Value = 6311
$ws.Range("D10").Value = 867961
$ws.Range("E10").Value = 68820
$ws.Range("G10").Value = 132
$ws.Range("H10").Value = 29102

# Row 15
$ws.Range("B15").Value = 705254
$ws.Range("C15").Value = 1461
$ws.Range("D15").Value = 635257
$ws.Range("E15").Value = 51505
$ws.Range("G15").Value = 21
$ws.Range("H15").Value = 18492

# Row 46
$ws.Range("B46").Value = 105547
$ws.Range("C46").Value = 123
$ws.Range("D46").Value = 98314
$ws.Range("E46").Value = 1103
$ws.Range("G46").Value = 10
$ws.Range("H46").Value = 6130

# Row 50
$ws.Range("B50").Value = 97075
$ws.Range("C50").Value = 632
$ws.Range("D50").Value = 59580
$ws.Range("E50").Value = 36291
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 1204

# Row 64
$ws.Range("B64").Value = 61558
$ws.Range("C64").Value = 118
$ws.Range("D64").Value = 56697
$ws.Range("E64").Value = 3736

# Row 76
$ws.Range("B76").Value = 42727
$ws.Range("C76").Value = 2185
$ws.Range("E76").Value = 37008
$ws.Range("G76").Value = 61
$ws.Range("H76").Value = 687

# Row 85
$ws.Range("B85").Value = 30527
$ws.Range("C85").Value = 1024
$ws.Range("D85").Value = 17153
$ws.Range("E85").Value = 12366
$ws.Range("G85").Value = 22
$ws.Range("H85").Value = 1008

# Row 97
$ws.Range("B97").Value = 16603
$ws.Range("C97").Value = 147
$ws.Range("E97").Value = 4462

# Row 146
$ws.Range("B146").Value = 3765
$ws.Range("C146").Value = 31
$ws.Range("D146").Value = 2749
$ws.Range("E146").Value = 905
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 111

# Row 151
$ws.Range("B151").Value = 3407
$ws.Range("C151").Value = 19
$ws.Range("D151").Value = 2588
$ws.Range("E151").Value = 687

# Row 158
$ws.Range("B158").Value = 2387
$ws.Range("C158").Value = 6
$ws.Range("D158").Value = 1802
$ws.Range("E158").Value = 520

# Row 161
$ws.Range("B161").Value = 2071
$ws.Range("C161").Value = 14
$ws.Range("D161").Value = 1541
$ws.Range("E161").Value = 479

# Row 167
$ws.Range("B167").Value = 1211
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 1128
$ws.Range("E167").Value = 14

# Row 175
$ws.Range("D175").Value = 541
$ws.Range("E175").Value = 33
"Update complete"
